$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for Price (D) cells whose new values would otherwise be
# auto-coerced to numbers (losing literal formatting like trailing zeros).
$ws.Range("D4:D7").NumberFormat = "@"
$ws.Range("D9:D14").NumberFormat = "@"
$ws.Range("D16:D17").NumberFormat = "@"
$ws.Range("D21:D26").NumberFormat = "@"
$ws.Range("D29:D37").NumberFormat = "@"
$ws.Range("D39:D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44:D45").NumberFormat = "@"
$ws.Range("D47:D51").NumberFormat = "@"

# Apply the updated coin data.
$ws.Range("D2").Value = "69.567.11"
$ws.Range("E2").Value = "  +1.06%  "

$ws.Range("D3").Value = "3.529.23"
$ws.Range("E3").Value = "  +0.93%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").Value = "606.06"
$ws.Range("E5").Value = "  +4.88%  "

$ws.Range("D6").Value = "171.83"
$ws.Range("E6").Value = "  -2.96%  "

$ws.Range("D7").Value = "0.615"
$ws.Range("E7").Value = "  -0.77%  "

$ws.Range("D8").Value = "3.525.80"
$ws.Range("E8").Value = "  +0.99%  "

$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("D10").Value = "0.196"
$ws.Range("E10").Value = "  +4.37%  "

$ws.Range("D11").Value = "6.69"
$ws.Range("E11").Value = "  +0.68%  "

$ws.Range("D12").Value = "0.583"
$ws.Range("E12").Value = "  -3.34%  "

$ws.Range("D13").Value = "47.58"
$ws.Range("E13").Value = "  +0.53%  "

$ws.Range("D14").Value = "0.0000279"
$ws.Range("E14").Value = "  +0.93%  "

$ws.Range("D15").Value = "4.100.42"
$ws.Range("E15").Value = "  +1.16%  "

$ws.Range("D16").Value = "8.40"
$ws.Range("E16").Value = "  -5.84%  "

$ws.Range("D17").Value = "617.91"
$ws.Range("E17").Value = "  -10.09%  "

$ws.Range("D18").Value = "3.534.17"
$ws.Range("E18").Value = "  +1.22%  "

$ws.Range("D19").Value = "69.806.82"
$ws.Range("E19").Value = "  +1.46%  "

$ws.Range("E20").Value = "  -1.36%  "

$ws.Range("D21").Value = "17.31"
$ws.Range("E21").Value = "  -1.22%  "

$ws.Range("D22").Value = "11.26"
$ws.Range("E22").Value = "  +0.83%  "

$ws.Range("D23").Value = "0.885"
$ws.Range("E23").Value = "  -2.08%  "

$ws.Range("D24").Value = "15.85"
$ws.Range("E24").Value = "  -3.04%  "

$ws.Range("D25").Value = "96.66"
$ws.Range("E25").Value = "  -1.37%  "

$ws.Range("D26").Value = "3.89"
$ws.Range("E26").Value = "  +1.34%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("E28").Value = "  -0.61%  "

$ws.Range("D29").Value = "9.31"
$ws.Range("E29").Value = "  -1.09%  "

$ws.Range("D30").Value = "33.34"
$ws.Range("E30").Value = "  +0.92%  "

$ws.Range("D31").Value = "3.14"
$ws.Range("E31").Value = "  -1.40%  "

$ws.Range("D32").Value = "8.51"
$ws.Range("E32").Value = "  -2.75%  "

$ws.Range("D33").Value = "1.34"
$ws.Range("E33").Value = "  -0.50%  "

$ws.Range("D34").Value = "6.98"
$ws.Range("E34").Value = "  -4.65%  "

$ws.Range("D35").Value = "571.14"
$ws.Range("E35").Value = "  +0.95%  "

$ws.Range("D36").Value = "3.59"
$ws.Range("E36").Value = "  -2.20%  "

$ws.Range("D37").Value = "10.82"
$ws.Range("E37").Value = "  -1.12%  "

$ws.Range("E38").Value = "  -3.20%  "

$ws.Range("D39").Value = "57.20"
$ws.Range("E39").Value = "  +0.56%  "

$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.29%  "

$ws.Range("E41").Value = "  +2.12%  "

$ws.Range("D42").Value = "0.0447"
$ws.Range("E42").Value = "  +1.54%  "

$ws.Range("D43").Value = "3.391.57"
$ws.Range("E43").Value = "  -1.07%  "

$ws.Range("D44").Value = "0.327"
$ws.Range("E44").Value = "  -2.65%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "33.18"
$ws.Range("E45").Value = "  -0.91%  "

$ws.Range("B46").Value = "PEPE"
$ws.Range("C46").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D46").Value = "0.0₃0709"
$ws.Range("E46").Value = "  +0.80%  "

$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").Value = "2.90"
$ws.Range("E47").Value = "  -1.39%  "

$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").Value = "2.60"
$ws.Range("E48").Value = "  -0.08%  "

$ws.Range("D49").Value = "0.129"
$ws.Range("E49").Value = "  -2.88%  "

$ws.Range("D50").Value = "134.21"
$ws.Range("E50").Value = "  -0.08%  "

$ws.Range("D51").Value = "5.68"
$ws.Range("E51").Value = "  +9.58%  "

